$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes current rows 3,4 down to 4,5)
$ws.Rows.Item(3).Insert()

# Copy the date cell style (column D uses a date-formatted style) from D4 (old D3, now shifted to D4) into new D3
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the newly inserted row 3 with the new weekly record
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44708
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100107
$ws.Cells.Item(3, 8).Value = "Otros"
$ws.Cells.Item(3, 9).Value = 100107001
$ws.Cells.Item(3, 10).Value = "Caqui"
$ws.Cells.Item(3, 11).Value = "Mankaki"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 70
$ws.Cells.Item(3, 14).Value = 12000
$ws.Cells.Item(3, 15).Value = 13000
$ws.Cells.Item(3, 16).Value = 12571
$ws.Cells.Item(3, 17).Value = "`$/caja 12 kilos empedrada"
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(3, 19).Value = 1048
$ws.Cells.Item(3, 20).Value = 12
